$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (NATMI LR-pairs recompute: Thbs1-Itgb1)

# Row 2
$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 3086.815927941086
$ws.Range("R2").Value = 27781.34335146977
$ws.Range("S2").Value = 0.02512035019736454
$ws.Range("T2").Value = 0.02512035019736454

# Row 3
$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 3576.190444429433
$ws.Range("R3").Value = 32185.7139998649
$ws.Range("S3").Value = 0.02910285499156939
$ws.Range("T3").Value = 0.02910285499156939

# Row 4
$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 2714.480670059282
$ws.Range("R4").Value = 24430.32603053354
$ws.Range("S4").Value = 0.02209030490565985
$ws.Range("T4").Value = 0.02209030490565985

# Row 5
$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1393.21419777364
$ws.Range("R5").Value = 12538.92777996276
$ws.Range("S5").Value = 0.01133790590855151
$ws.Range("T5").Value = 0.01133790590855151

# Row 6
$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 17749.39472892449
$ws.Range("R6").Value = 159744.5525603205
$ws.Range("S6").Value = 0.1444436668043361
$ws.Range("T6").Value = 0.1444436668043361

# Row 7
$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.1673433314595345
$ws.Range("T7").Value = 0.1673433314595345

# Row 8
$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 15608.44249273202
$ws.Range("R8").Value = 140475.9824345881
$ws.Range("S8").Value = 0.1270207069698447
$ws.Range("T8").Value = 0.1270207069698447

# Row 9
$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 8011.073324582828
$ws.Range("R9").Value = 72099.65992124545
$ws.Range("S9").Value = 0.06519370512141758
$ws.Range("T9").Value = 0.06519370512141758

# Row 10
$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 5421.11417291511
$ws.Range("R10").Value = 48790.027556236
$ws.Range("S10").Value = 0.04411674996583176
$ws.Range("T10").Value = 0.04411674996583177

# Row 11
$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 6280.561314931119
$ws.Range("R11").Value = 56525.05183438008
$ws.Range("S11").Value = 0.05111088686532825
$ws.Range("T11").Value = 0.05111088686532826

# Row 12
$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 4767.213198351534
$ws.Range("R12").Value = 42904.91878516381
$ws.Range("S12").Value = 0.03879533726779915
$ws.Range("T12").Value = 0.03879533726779916

# Row 13
$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 2446.784456789737
$ws.Range("R13").Value = 22021.06011110764
$ws.Range("S13").Value = 0.01991180680897396
$ws.Range("T13").Value = 0.01991180680897396

# Row 14
$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 8959.629655362116
$ws.Range("R14").Value = 80636.66689825903
$ws.Range("S14").Value = 0.07291300804305186
$ws.Range("T14").Value = 0.07291300804305186

# Row 15
$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 10380.06240317162
$ws.Range("R15").Value = 93420.56162854456
$ws.Range("S15").Value = 0.08447241712015199
$ws.Range("T15").Value = 0.08447241712015199

# Row 16
$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 7878.908907468405
$ws.Range("R16").Value = 70910.18016721564
$ws.Range("S16").Value = 0.06411815785230673
$ws.Range("T16").Value = 0.06411815785230673

# Row 17
$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 4043.87029678515
$ws.Range("R17").Value = 24430.32603053354
$ws.Range("S17").Value = 0.03290880971827818
$ws.Range("T17").Value = 0.03290880971827818
